$wb = $excel.ActiveWorkbook

# --- Schedule sheet updates ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1043.1788055
$wsSchedule.Range("F2").Value = 22.99776908068783
$wsSchedule.Range("E3").Value = 368.2275285
$wsSchedule.Range("F3").Value = 24.35367251984127
$wsSchedule.Range("E4").Value = -58.6619475
$wsSchedule.Range("F4").Value = -1.293252810846561

# --- Detailed sheet updates ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B29").Value = 36.0601
$wsDetailed.Range("B30").Value = 36.0601
$wsDetailed.Range("B31").Value = -4.05321
$wsDetailed.Range("C31").Value = "historical"
$wsDetailed.Range("B32").Value = -15.51759
$wsDetailed.Range("C32").Value = "historical"
$wsDetailed.Range("B33").Value = -6.57876
$wsDetailed.Range("C33").Value = "historical"
$wsDetailed.Range("B34").Value = 4.29232
$wsDetailed.Range("B35").Value = -3.47865
$wsDetailed.Range("B36").Value = -6
$wsDetailed.Range("B37").Value = -1.57063
$wsDetailed.Range("B38").Value = -3.20521
$wsDetailed.Range("B39").Value = 1.18622
$wsDetailed.Range("B40").Value = 0.00045
$wsDetailed.Range("B41").Value = 23.74544
$wsDetailed.Range("B42").Value = 32.79309
$wsDetailed.Range("B43").Value = 24.74099
$wsDetailed.Range("B44").Value = 21.88816
$wsDetailed.Range("B46").Value = 64.8901
$wsDetailed.Range("B47").Value = 62.33685
$wsDetailed.Range("B61").Value = 57.06017
$wsDetailed.Range("B62").Value = 51.15669
$wsDetailed.Range("B64").Value = 36.06029
$wsDetailed.Range("B65").Value = 6.27504
$wsDetailed.Range("B66").Value = -4.56332
$wsDetailed.Range("B67").Value = -0.90384
$wsDetailed.Range("B68").Value = 0.0094
$wsDetailed.Range("B69").Value = -0.00776
$wsDetailed.Range("B70").Value = -0.89434
$wsDetailed.Range("B71").Value = -5.50985
$wsDetailed.Range("B73").Value = 6.25571
$wsDetailed.Range("B74").Value = 8.71008
$wsDetailed.Range("B75").Value = 0.7
$wsDetailed.Range("B76").Value = -5.51011
$wsDetailed.Range("B77").Value = -7.01
$wsDetailed.Range("B78").Value = -11.17408
$wsDetailed.Range("B79").Value = -12.11173
$wsDetailed.Range("B80").Value = -6.72418
$wsDetailed.Range("B81").Value = -5.66385
$wsDetailed.Range("B82").Value = -2.21718
$wsDetailed.Range("B83").Value = -5.17419
$wsDetailed.Range("B84").Value = -6.49855
$wsDetailed.Range("B85").Value = -3.05311
$wsDetailed.Range("B86").Value = -2.93302
$wsDetailed.Range("B87").Value = -2.92219
$wsDetailed.Range("B88").Value = 0.75497
$wsDetailed.Range("B89").Value = 21.60312
$wsDetailed.Range("B90").Value = 51.84959
$wsDetailed.Range("B91").Value = 32.62903
$wsDetailed.Range("B92").Value = 29.54841
$wsDetailed.Range("B93").Value = 36.0601
$wsDetailed.Range("B94").Value = 36.06045
$wsDetailed.Range("B95").Value = 57.3
$wsDetailed.Range("B96").Value = 57.3
$wsDetailed.Range("B97").Value = 57.06005
